# Solid state eem correction performed
# Rename BBWMO1181-series sample names to BBWMO181-series (drop the extra "1")
# throughout the log sheet, and fix the dependent "Corrected Name" formula
# text in C33 that referenced the old naming scheme.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: "Name of Raw EEM" corrections (BBWMO1181... -> BBWMO181...)
$ws.Range("A5").Value  = "BBWMO181"
$ws.Range("A6").Value  = "BBWMO18110X"
$ws.Range("A7").Value  = "BBWMO18120X"
$ws.Range("A8").Value  = "BBWMO18140X"

$ws.Range("A12").Value = "BBWMO181"
$ws.Range("A13").Value = "BBWMO18110X"
$ws.Range("A14").Value = "BBWMO18120X"
$ws.Range("A15").Value = "BBWMO18140X"

$ws.Range("A19").Value = "BBWMO181"
$ws.Range("A20").Value = "BBWMO18110X"
$ws.Range("A21").Value = "BBWMO18120X"
$ws.Range("A22").Value = "BBWMO18140X"

$ws.Range("A26").Value = "BBWMO181"
$ws.Range("A27").Value = "BBWMO18110X"
$ws.Range("A28").Value = "BBWMO18120X"
$ws.Range("A29").Value = "BBWMO18140X"

$ws.Range("A33").Value = "BBWMO181"
$ws.Range("A34").Value = "BBWMO18110X"
$ws.Range("A35").Value = "BBWMO18120X"
$ws.Range("A36").Value = "BBWMO18140X"

$ws.Range("A40").Value = "BBWMO181"
$ws.Range("A41").Value = "BBWMO18110X"
$ws.Range("A42").Value = "BBWMO18120X"
$ws.Range("A43").Value = "BBWMO18140X"

# Column C: "Corrected Name" dependent value that also needed updating
$ws.Range("C33").Value = "COMBSMPLHLDRBBWM0181"

# Restore view/selection state captured in the saved workbook
$ws.Range("C30").Select()
